$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string for tc12 -> tc12_queue across all cells that use it (D13:D15)
$ws.Range("D13:D15").Value = "Implementation tc12_queue"

# Update row 15 effort numbers
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 3

# Add new row 16
$ws.Range("A16").Value = 41444
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").Value = 1.5
$ws.Range("C16").Value = 2.5
$ws.Range("D16").Value = "Design and implementation tc13_eventStates"

# Update selection to match post-edit state (next empty row)
$ws.Range("A17").Select()
